# Test case added for refactored code:
# - Update the FX date in A2 (43951 -> 43950, i.e. 2020-04-30 -> 2020-04-29)
# - Move the active selection to A3

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("A2").Value = 43950

$ws.Range("A3").Select()
